$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$shp = $s.Shapes.AddTextbox(1, 605.34, 575.92, 754.04, 35.34)
$shp.Name = "TextBox 5"
$shp.TextFrame.MarginLeft = 0
$shp.TextFrame.MarginRight = 0
$shp.TextFrame.MarginTop = 0
$shp.TextFrame.MarginBottom = 0
$shp.TextFrame.TextRange.Text = "{{ tendering }}"
$shp.TextFrame.TextRange.Font.Name = "Times New Roman"
$shp.TextFrame.TextRange.Font.Size = 32
$shp.TextFrame.AutoSize = 1
